# Add three newly-collected submissions to the raw data sheet
# ("八位序列号收集收集结果yd5", the first/raw sheet of the workbook).
# Columns: A = submitter, B = submit time, C = serial number, D = QQ number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 56: 文昊 / 883c8e20 / 1243776890 -------------------------------
$ws.Cells.Item(56, 1).Value = "文昊"

$ws.Cells.Item(56, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(56, 2).Value = 45931.129525463

$ws.Cells.Item(56, 3).Value = "883c8e20"

$ws.Cells.Item(56, 4).NumberFormat = "@"
$ws.Cells.Item(56, 4).Value = "1243776890"

# --- Row 57: clh / 6cffec45 / 1007021745 --------------------------------
$ws.Cells.Item(57, 1).Value = "clh"

$ws.Cells.Item(57, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(57, 2).Value = 45931.5039236111

$ws.Cells.Item(57, 3).Value = "6cffec45"

$ws.Cells.Item(57, 4).NumberFormat = "@"
$ws.Cells.Item(57, 4).Value = "1007021745"

# --- Row 58: Non-numb / f8241963 / 1183413694 ---------------------------
$ws.Cells.Item(58, 1).Value = "Non-numb"

$ws.Cells.Item(58, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(58, 2).Value = 45931.6418634259

$ws.Cells.Item(58, 3).Value = "f8241963"

$ws.Cells.Item(58, 4).NumberFormat = "@"
$ws.Cells.Item(58, 4).Value = "1183413694"
